$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 16410575
$ws.Range("I15").Value = 16410575
$ws.Range("K15").Value = 49231725
$ws.Range("M15").Value = -49231556
$ws.Range("H28").Value = 1377.5
$ws.Range("J28").Value = 250
$ws.Range("L28").Value = 250
$ws.Range("N28").Value = -1220
$ws.Range("H32").Value = 1651.381
$ws.Range("I32").Value = 1873.75
$ws.Range("J32").Value = 939.8
$ws.Range("K32").Value = 1873.75
$ws.Range("L32").Value = 939.8
$ws.Range("M32").Value = -1547.75
$ws.Range("N32").Value = -1591.8
$ws.Range("H38").Value = 161.8
$ws.Range("J38").Value = 97
$ws.Range("L38").Value = 291
$ws.Range("N38").Value = -1035
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("N69").ClearContents()
$ws.Range("H70").Value = 53033348
$ws.Range("J70").Value = 66672350
$ws.Range("L70").Value = 200017050
$ws.Range("N70").Value = -200017590
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("N72").ClearContents()
$ws.Range("H73").Value = 53033348
$ws.Range("J73").Value = 66672350
$ws.Range("L73").Value = 200017050
$ws.Range("N73").Value = -200018922
$ws.Range("H137").Value = 4291.5264
$ws.Range("I137").Value = 5774.6665
$ws.Range("J137").Value = 2956.7
$ws.Range("K137").Value = 17323.9995
$ws.Range("L137").Value = 8870.099999999999
$ws.Range("M137").Value = -14773.9995
$ws.Range("N137").Value = -13970.1
$ws.Range("H138").Value = 1497801.8
$ws.Range("I138").Value = 3128.3333
$ws.Range("J138").Value = 2046865.5
$ws.Range("K138").Value = 9384.999899999999
$ws.Range("L138").Value = 6140596.5
$ws.Range("M138").Value = -4244.999899999999
$ws.Range("N138").Value = -6150876.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3092.077
$ws.Range("J32").Value = 4125
$ws.Range("L32").Value = 4125
$ws.Range("N32").Value = -4699
$ws.Range("H61").Value = 4443.013
$ws.Range("I61").Value = 2563.2131
$ws.Range("K61").Value = 2563.2131
$ws.Range("M61").Value = -2351.2131
$ws.Range("H110").Value = 25644260
$ws.Range("I110").Value = 3291.7144
$ws.Range("K110").Value = 3291.7144
$ws.Range("M110").Value = -1246.7144
$ws.Range("H132").Value = 765884.4
$ws.Range("I132").Value = 1145352.8
$ws.Range("K132").Value = 3436058.4
$ws.Range("M132").Value = -3433528.4
$ws.Range("H136").Value = 4443.013
$ws.Range("I136").Value = 2563.2131
$ws.Range("K136").Value = 7689.6393
$ws.Range("M136").Value = -5139.6393

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H24").Value = 2001
$ws.Range("I24").Value = 2001
$ws.Range("K24").Value = 2001
$ws.Range("M24").Value = -1766
$ws.Range("H29").Value = 25322.5
$ws.Range("I29").Value = 430.33334
$ws.Range("K29").Value = 430.33334
$ws.Range("M29").Value = -141.33334
$ws.Range("H64").Value = 23810640
$ws.Range("I64").Value = 47619844
$ws.Range("K64").Value = 47619844
$ws.Range("M64").Value = -47619619
$ws.Range("H67").Value = 23810640
$ws.Range("I67").Value = 47619844
$ws.Range("K67").Value = 47619844
$ws.Range("M67").Value = -47619064
$ws.Range("H105").Value = 3799.6667
$ws.Range("I105").Value = 1514.3334
$ws.Range("K105").Value = 1514.3334
$ws.Range("M105").Value = 232.6666
$ws.Range("H134").Value = 3794.578
$ws.Range("I134").Value = 1194.875
$ws.Range("K134").Value = 3584.625
$ws.Range("M134").Value = -1049.625

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 374.16666
$ws.Range("I22").Value = 296
$ws.Range("K22").Value = 296
$ws.Range("M22").Value = 54
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H99").Value = 8200.546
$ws.Range("I99").Value = 10383.167
$ws.Range("K99").Value = 10383.167
$ws.Range("M99").Value = -8885.166999999999
$ws.Range("H126").Value = 8200.546
$ws.Range("I126").Value = 10383.167
$ws.Range("K126").Value = 31149.501
$ws.Range("M126").Value = -28679.501
$ws.Range("H132").Value = 3068.6892
$ws.Range("I132").Value = 1699.386
$ws.Range("J132").Value = 7659.8823
$ws.Range("K132").Value = 5098.157999999999
$ws.Range("L132").Value = 22979.6469
$ws.Range("M132").Value = -2568.157999999999
$ws.Range("N132").Value = -28039.6469
$ws.Range("H134").Value = 4559.6113
$ws.Range("I134").Value = 1901.9722
$ws.Range("K134").Value = 5705.9166
$ws.Range("M134").Value = -3170.9166

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 83085656
$ws.Range("I4").Value = 76456264
$ws.Range("K4").Value = 229368792
$ws.Range("M4").Value = -229368680
$ws.Range("H8").Value = 714.4286
$ws.Range("I8").Value = 714.4286
$ws.Range("K8").Value = 2143.2858
$ws.Range("M8").Value = -2004.2858
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H129").Value = 106152190
$ws.Range("I129").Value = 827
$ws.Range("K129").Value = 2481
$ws.Range("M129").Value = 2519
$ws.Range("H130").Value = 0
$ws.Range("I130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("M130").ClearContents()
$ws.Range("H138").Value = 255982.25
$ws.Range("I138").Value = 255982.25
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 767946.75
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = -762806.75
$ws.Range("N138").ClearContents()
$ws.Range("H139").Value = 162178.16
$ws.Range("I139").Value = 203092.4
$ws.Range("K139").Value = 609277.2
$ws.Range("M139").Value = -604137.2

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 2500150
$ws.Range("I2").Value = 120.4
$ws.Range("K2").Value = 120.4
$ws.Range("M2").Value = -7.400000000000006
$ws.Range("H80").Value = 2430.3845
$ws.Range("I80").Value = 2731.75
$ws.Range("K80").Value = 2731.75
$ws.Range("M80").Value = -1733.75
$ws.Range("H83").Value = 2430.3845
$ws.Range("I83").Value = 2731.75
$ws.Range("K83").Value = 13658.75
$ws.Range("M83").Value = -8666.75
$ws.Range("H132").Value = 4507.472
$ws.Range("I132").Value = 3347.652
$ws.Range("J132").Value = 6559.4614
$ws.Range("K132").Value = 10042.956
$ws.Range("L132").Value = 19678.3842
$ws.Range("M132").Value = -7512.956
$ws.Range("N132").Value = -24738.3842

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3024.9375
$ws.Range("I22").Value = 1881.1428
$ws.Range("J22").Value = 3914.5557
$ws.Range("K22").Value = 1881.1428
$ws.Range("L22").Value = 3914.5557
$ws.Range("M22").Value = -1586.1428
$ws.Range("N22").Value = -4504.5557
$ws.Range("H27").Value = 3024.9375
$ws.Range("I27").Value = 1881.1428
$ws.Range("J27").Value = 3914.5557
$ws.Range("K27").Value = 1881.1428
$ws.Range("L27").Value = 3914.5557
$ws.Range("M27").Value = -1774.1428
$ws.Range("N27").Value = -4128.5557

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5229.1665
$ws.Range("I62").Value = 4535
$ws.Range("J62").Value = 6201
$ws.Range("K62").Value = 4535
$ws.Range("L62").Value = 6201
$ws.Range("M62").Value = -3911
$ws.Range("N62").Value = -7449
$ws.Range("H65").Value = 5229.1665
$ws.Range("I65").Value = 4535
$ws.Range("J65").Value = 6201
$ws.Range("K65").Value = 22675
$ws.Range("L65").Value = 31005
$ws.Range("M65").Value = -19555
$ws.Range("N65").Value = -37245
$ws.Range("H132").Value = 7476679
$ws.Range("I132").Value = 9266475
$ws.Range("K132").Value = 27799425
$ws.Range("M132").Value = -27796895
